$wb = $excel.ActiveWorkbook

# --- Update the conversion summary text on "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.01 = 15449.44 pesos`n✅ 15449.44 pesos = 4.0 = 960.49 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate values on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 249.2
$wsTasas.Range("O10").Value = 3850
$wsTasas.Range("N12").Value = 3862
$wsTasas.Range("O12").Value = 240.1
